$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(58, 8).Value2 = 50
$ws.Cells.Item(58, 10).Value2 = 0
$ws.Cells.Item(58, 12).Value2 = 0
$ws.Cells.Item(58, 14).Value2 = $null
$ws.Cells.Item(64, 8).Value2 = 1500
$ws.Cells.Item(64, 10).Value2 = 1500
$ws.Cells.Item(64, 12).Value2 = 1500
$ws.Cells.Item(64, 14).Value2 = -1996
$ws.Cells.Item(67, 8).Value2 = 1500
$ws.Cells.Item(67, 10).Value2 = 1500
$ws.Cells.Item(67, 12).Value2 = 1500
$ws.Cells.Item(67, 14).Value2 = -3216
$ws.Cells.Item(95, 8).Value2 = 29378.8
$ws.Cells.Item(95, 10).Value2 = 29378.8
$ws.Cells.Item(95, 12).Value2 = 29378.8
$ws.Cells.Item(95, 14).Value2 = -34870.8
$ws.Cells.Item(116, 8).Value2 = 7711.4287
$ws.Cells.Item(116, 10).Value2 = 9996
$ws.Cells.Item(116, 12).Value2 = 9996
$ws.Cells.Item(116, 14).Value2 = -16880
$ws.Cells.Item(124, 8).Value2 = 0
$ws.Cells.Item(124, 10).Value2 = 0
$ws.Cells.Item(124, 12).Value2 = $null
$ws.Cells.Item(124, 14).Value2 = 0
$ws.Cells.Item(137, 8).Value2 = 6862.909
$ws.Cells.Item(137, 9).Value2 = 3420
$ws.Cells.Item(137, 11).Value2 = 10260
$ws.Cells.Item(137, 13).Value2 = -7710

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(38, 8).Value2 = 496.5
$ws.Cells.Item(38, 10).Value2 = 496.5
$ws.Cells.Item(38, 12).Value2 = 496.5
$ws.Cells.Item(38, 14).Value2 = -1430.5
$ws.Cells.Item(39, 8).Value2 = 6500
$ws.Cells.Item(39, 10).Value2 = 0
$ws.Cells.Item(39, 12).Value2 = 0
$ws.Cells.Item(39, 14).Value2 = $null
$ws.Cells.Item(41, 8).Value2 = 2432.5
$ws.Cells.Item(41, 10).Value2 = 3599
$ws.Cells.Item(41, 12).Value2 = 3599
$ws.Cells.Item(41, 14).Value2 = -4427
$ws.Cells.Item(42, 8).Value2 = 0
$ws.Cells.Item(42, 10).Value2 = 0
$ws.Cells.Item(42, 12).Value2 = $null
$ws.Cells.Item(42, 14).Value2 = 0
$ws.Cells.Item(61, 8).Value2 = 8624.25
$ws.Cells.Item(61, 9).Value2 = 5832.3335
$ws.Cells.Item(61, 11).Value2 = 5832.3335
$ws.Cells.Item(61, 13).Value2 = -5620.3335
$ws.Cells.Item(74, 8).Value2 = 2998.75
$ws.Cells.Item(74, 9).Value2 = 2331.6667
$ws.Cells.Item(74, 10).Value2 = 5000
$ws.Cells.Item(74, 11).Value2 = 2331.6667
$ws.Cells.Item(74, 12).Value2 = 5000
$ws.Cells.Item(74, 13).Value2 = -1457.6667
$ws.Cells.Item(74, 14).Value2 = -6748
$ws.Cells.Item(77, 8).Value2 = 2998.75
$ws.Cells.Item(77, 9).Value2 = 2331.6667
$ws.Cells.Item(77, 10).Value2 = 5000
$ws.Cells.Item(77, 11).Value2 = 11658.3335
$ws.Cells.Item(77, 12).Value2 = 25000
$ws.Cells.Item(77, 13).Value2 = -7290.333500000001
$ws.Cells.Item(77, 14).Value2 = -33736
$ws.Cells.Item(112, 8).Value2 = 30000
$ws.Cells.Item(112, 10).Value2 = 30000
$ws.Cells.Item(112, 12).Value2 = 30000
$ws.Cells.Item(112, 14).Value2 = -32954
$ws.Cells.Item(132, 8).Value2 = 15306.909
$ws.Cells.Item(132, 9).Value2 = 13482.286
$ws.Cells.Item(132, 11).Value2 = 40446.858
$ws.Cells.Item(132, 13).Value2 = -37916.858
$ws.Cells.Item(136, 8).Value2 = 8624.25
$ws.Cells.Item(136, 9).Value2 = 5832.3335
$ws.Cells.Item(136, 11).Value2 = 17497.0005
$ws.Cells.Item(136, 13).Value2 = -14947.0005

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(41, 8).Value2 = 0
$ws.Cells.Item(41, 10).Value2 = 0
$ws.Cells.Item(41, 12).Value2 = $null
$ws.Cells.Item(41, 14).Value2 = 0
$ws.Cells.Item(48, 8).Value2 = 0
$ws.Cells.Item(48, 10).Value2 = 0
$ws.Cells.Item(48, 12).Value2 = $null
$ws.Cells.Item(48, 14).Value2 = 0
$ws.Cells.Item(86, 8).Value2 = 0
$ws.Cells.Item(86, 9).Value2 = 0
$ws.Cells.Item(86, 11).Value2 = 0
$ws.Cells.Item(86, 13).Value2 = $null
$ws.Cells.Item(89, 8).Value2 = 0
$ws.Cells.Item(89, 9).Value2 = 0
$ws.Cells.Item(89, 11).Value2 = 0
$ws.Cells.Item(89, 13).Value2 = $null
$ws.Cells.Item(94, 8).Value2 = 0
$ws.Cells.Item(94, 9).Value2 = 0
$ws.Cells.Item(94, 11).Value2 = 0
$ws.Cells.Item(94, 13).Value2 = $null
$ws.Cells.Item(134, 8).Value2 = 10396.714
$ws.Cells.Item(134, 9).Value2 = 4444.25
$ws.Cells.Item(134, 10).Value2 = 18333.334
$ws.Cells.Item(134, 11).Value2 = 13332.75
$ws.Cells.Item(134, 12).Value2 = 55000.00199999999
$ws.Cells.Item(134, 13).Value2 = -10797.75
$ws.Cells.Item(134, 14).Value2 = -60070.00199999999

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value2 = 5159.9375
$ws.Cells.Item(31, 10).Value2 = 8256.333000000001
$ws.Cells.Item(31, 12).Value2 = 8256.333000000001
$ws.Cells.Item(31, 14).Value2 = -8846.333000000001
$ws.Cells.Item(34, 8).Value2 = 5159.9375
$ws.Cells.Item(34, 10).Value2 = 8256.333000000001
$ws.Cells.Item(34, 12).Value2 = 8256.333000000001
$ws.Cells.Item(34, 14).Value2 = -8660.333000000001
$ws.Cells.Item(58, 8).Value2 = 14399.2
$ws.Cells.Item(58, 9).Value2 = 9998.5
$ws.Cells.Item(58, 11).Value2 = 9998.5
$ws.Cells.Item(58, 13).Value2 = -9795.5
$ws.Cells.Item(136, 8).Value2 = 14399.2
$ws.Cells.Item(136, 9).Value2 = 9998.5
$ws.Cells.Item(136, 11).Value2 = 29995.5
$ws.Cells.Item(136, 13).Value2 = -27445.5

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value2 = 42921324
$ws.Cells.Item(4, 9).Value2 = 290.2
$ws.Cells.Item(4, 11).Value2 = 870.5999999999999
$ws.Cells.Item(4, 13).Value2 = -758.5999999999999
$ws.Cells.Item(32, 8).Value2 = 214.28572
$ws.Cells.Item(32, 9).Value2 = 900
$ws.Cells.Item(32, 10).Value2 = 100
$ws.Cells.Item(32, 11).Value2 = 2700
$ws.Cells.Item(32, 12).Value2 = 300
$ws.Cells.Item(32, 13).Value2 = -2417
$ws.Cells.Item(32, 14).Value2 = -866
$ws.Cells.Item(48, 8).Value2 = 225
$ws.Cells.Item(48, 9).Value2 = 225
$ws.Cells.Item(48, 11).Value2 = 675
$ws.Cells.Item(48, 13).Value2 = -425
$ws.Cells.Item(68, 8).Value2 = 0
$ws.Cells.Item(68, 9).Value2 = 0
$ws.Cells.Item(68, 10).Value2 = 0
$ws.Cells.Item(68, 11).Value2 = 0
$ws.Cells.Item(68, 12).Value2 = $null
$ws.Cells.Item(68, 13).Value2 = $null
$ws.Cells.Item(68, 14).Value2 = 0
$ws.Cells.Item(71, 8).Value2 = 0
$ws.Cells.Item(71, 9).Value2 = 0
$ws.Cells.Item(71, 10).Value2 = 0
$ws.Cells.Item(71, 11).Value2 = 0
$ws.Cells.Item(71, 12).Value2 = $null
$ws.Cells.Item(71, 13).Value2 = $null
$ws.Cells.Item(71, 14).Value2 = 0
$ws.Cells.Item(98, 8).Value2 = 201.5
$ws.Cells.Item(98, 9).Value2 = 102
$ws.Cells.Item(98, 11).Value2 = 306
$ws.Cells.Item(98, 13).Value2 = 1192

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(44, 8).Value2 = 0
$ws.Cells.Item(44, 10).Value2 = 0
$ws.Cells.Item(44, 12).Value2 = $null
$ws.Cells.Item(44, 14).Value2 = 0
$ws.Cells.Item(126, 8).Value2 = 1998.5
$ws.Cells.Item(126, 9).Value2 = 1998.5
$ws.Cells.Item(126, 11).Value2 = 5995.5
$ws.Cells.Item(126, 13).Value2 = -3525.5

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value2 = 0
$ws.Cells.Item(46, 9).Value2 = 0
$ws.Cells.Item(46, 11).Value2 = 0
$ws.Cells.Item(46, 13).Value2 = $null
$ws.Cells.Item(68, 8).Value2 = 0
$ws.Cells.Item(68, 9).Value2 = 0
$ws.Cells.Item(68, 11).Value2 = 0
$ws.Cells.Item(68, 13).Value2 = $null
$ws.Cells.Item(71, 8).Value2 = 0
$ws.Cells.Item(71, 9).Value2 = 0
$ws.Cells.Item(71, 11).Value2 = 0
$ws.Cells.Item(71, 13).Value2 = $null
$ws.Cells.Item(97, 8).Value2 = 44163
$ws.Cells.Item(97, 10).Value2 = 44163
$ws.Cells.Item(97, 12).Value2 = 44163
$ws.Cells.Item(97, 14).Value2 = -46145
$ws.Cells.Item(100, 8).Value2 = 2000
$ws.Cells.Item(100, 9).Value2 = 2000
$ws.Cells.Item(100, 11).Value2 = 2000
$ws.Cells.Item(100, 13).Value2 = -1459

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value2 = 0
$ws.Cells.Item(62, 10).Value2 = 0
$ws.Cells.Item(62, 12).Value2 = $null
$ws.Cells.Item(62, 14).Value2 = 0
$ws.Cells.Item(65, 8).Value2 = 0
$ws.Cells.Item(65, 10).Value2 = 0
$ws.Cells.Item(65, 12).Value2 = $null
$ws.Cells.Item(65, 14).Value2 = 0
$ws.Cells.Item(132, 8).Value2 = 9332.777
$ws.Cells.Item(132, 9).Value2 = 7332.5
$ws.Cells.Item(132, 11).Value2 = 21997.5
$ws.Cells.Item(132, 13).Value2 = -19467.5
